$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6916.75
$ws.Range("J40").Value = 8180.3335
$ws.Range("L40").Value = 8180.3335
$ws.Range("N40").Value = -8530.333500000001
$ws.Range("H53").Value = 241.16667
$ws.Range("I53").Value = 234
$ws.Range("K53").Value = 234
$ws.Range("M53").Value = 403
$ws.Range("H76").Value = 4499.75
$ws.Range("I76").Value = 4333
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 4333
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -4018
$ws.Range("N76").Value = -5630
$ws.Range("H79").Value = 4499.75
$ws.Range("I79").Value = 4333
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 4333
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -3241
$ws.Range("N79").Value = -7184
$ws.Range("H86").Value = 3450
$ws.Range("I86").Value = 3450
$ws.Range("K86").Value = 3450
$ws.Range("M86").Value = -2327
$ws.Range("H89").Value = 3450
$ws.Range("I89").Value = 3450
$ws.Range("K89").Value = 17250
$ws.Range("M89").Value = -11634
$ws.Range("H100").Value = 2361.889
$ws.Range("I100").Value = 1969.625
$ws.Range("K100").Value = 1969.625
$ws.Range("M100").Value = -1428.625
$ws.Range("H103").Value = 791
$ws.Range("I103").Value = 785
$ws.Range("K103").Value = 2355
$ws.Range("M103").Value = -1769
$ws.Range("H138").Value = 5000
$ws.Range("J138").Value = 5000
$ws.Range("L138").Value = 15000
$ws.Range("N138").Value = -25280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5000
$ws.Range("I63").Value = 5000
$ws.Range("K63").Value = 5000
$ws.Range("M63").Value = -4314
$ws.Range("H66").Value = 5000
$ws.Range("I66").Value = 5000
$ws.Range("K66").Value = 25000
$ws.Range("M66").Value = -21568
$ws.Range("H97").Value = 913.75
$ws.Range("I97").Value = 615.7143
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 615.7143
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -119.7143
$ws.Range("N97").Value = -3992
$ws.Range("H102").Value = 1925
$ws.Range("I102").Value = 1925
$ws.Range("K102").Value = 1925
$ws.Range("M102").Value = -303
$ws.Range("H110").Value = 841.5
$ws.Range("I110").Value = 872.5454999999999
$ws.Range("J110").Value = 500
$ws.Range("K110").Value = 872.5454999999999
$ws.Range("L110").Value = 500
$ws.Range("M110").Value = 1172.4545
$ws.Range("N110").Value = -4590
$ws.Range("H137").Value = 49998
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3146.6
$ws.Range("I20").Value = 3146.6
$ws.Range("K20").Value = 3146.6
$ws.Range("M20").Value = -2899.6
$ws.Range("H26").Value = 15001
$ws.Range("I26").Value = 15001
$ws.Range("K26").Value = 15001
$ws.Range("M26").Value = -14709
$ws.Range("H80").Value = 395.33334
$ws.Range("I80").Value = 309.2
$ws.Range("K80").Value = 309.2
$ws.Range("M80").Value = 688.8
$ws.Range("H83").Value = 395.33334
$ws.Range("I83").Value = 309.2
$ws.Range("K83").Value = 1546
$ws.Range("M83").Value = 3446
$ws.Range("H94").Value = 4168.0625
$ws.Range("I94").Value = 2343.6365
$ws.Range("J94").Value = 8181.8
$ws.Range("K94").Value = 2343.6365
$ws.Range("L94").Value = 8181.8
$ws.Range("M94").Value = -1892.6365
$ws.Range("N94").Value = -9083.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 384.6316
$ws.Range("I7").Value = 573.7778
$ws.Range("J7").Value = 214.4
$ws.Range("K7").Value = 573.7778
$ws.Range("L7").Value = 214.4
$ws.Range("M7").Value = -460.7778
$ws.Range("N7").Value = -440.4
$ws.Range("H41").Value = 6000
$ws.Range("I41").Value = 6000
$ws.Range("K41").Value = 6000
$ws.Range("M41").Value = -5572
$ws.Range("H60").Value = 12717.647
$ws.Range("I60").Value = 9875
$ws.Range("K60").Value = 9875
$ws.Range("M60").Value = -9364
$ws.Range("H62").Value = 35720716
$ws.Range("I62").Value = 50007200
$ws.Range("K62").Value = 50007200
$ws.Range("M62").Value = -50006576
$ws.Range("H65").Value = 35720716
$ws.Range("I65").Value = 50007200
$ws.Range("K65").Value = 250036000
$ws.Range("M65").Value = -250032880
$ws.Range("H132").Value = 5165
$ws.Range("I132").Value = 4899
$ws.Range("K132").Value = 14697
$ws.Range("M132").Value = -12167
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 32.128204
$ws.Range("J2").Value = 61.923077
$ws.Range("L2").Value = 371.538462
$ws.Range("N2").Value = -597.538462
$ws.Range("H68").Value = 2292
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2292
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 6876
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -8498
$ws.Range("H71").Value = 2292
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2292
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 20628
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -28740
$ws.Range("H80").Value = 5454.364
$ws.Range("I80").Value = 3999.3333
$ws.Range("K80").Value = 11997.9999
$ws.Range("M80").Value = -11061.9999
$ws.Range("H83").Value = 5454.364
$ws.Range("I83").Value = 3999.3333
$ws.Range("K83").Value = 35993.9997
$ws.Range("M83").Value = -31313.9997
$ws.Range("H92").Value = 2300
$ws.Range("I92").Value = 2300
$ws.Range("K92").Value = 6900
$ws.Range("M92").Value = -5652

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 102003
$ws.Range("I80").Value = 4000
$ws.Range("K80").Value = 4000
$ws.Range("M80").Value = -3002
$ws.Range("H83").Value = 102003
$ws.Range("I83").Value = 4000
$ws.Range("K83").Value = 20000
$ws.Range("M83").Value = -15008
$ws.Range("H97").Value = 1852.9
$ws.Range("I97").Value = 1852.9
$ws.Range("K97").Value = 1852.9
$ws.Range("M97").Value = -1356.9
$ws.Range("H107").Value = 2433.3333
$ws.Range("I107").Value = 300
$ws.Range("K107").Value = 300
$ws.Range("M107").Value = 1620

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 2500
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = -2205
$ws.Range("H27").Value = 2500
$ws.Range("I27").Value = 2500
$ws.Range("K27").Value = 2500
$ws.Range("M27").Value = -2393
$ws.Range("H46").Value = 4046.182
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4046.182
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4046.182
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4422.182
$ws.Range("H55").Value = 3186.125
$ws.Range("I55").Value = 2999
$ws.Range("K55").Value = 2999
$ws.Range("M55").Value = -2826
$ws.Range("H68").Value = 4636.1816
$ws.Range("I68").Value = 2749.75
$ws.Range("J68").Value = 9666.666999999999
$ws.Range("K68").Value = 2749.75
$ws.Range("L68").Value = 9666.666999999999
$ws.Range("M68").Value = -2000.75
$ws.Range("N68").Value = -11164.667
$ws.Range("H71").Value = 4636.1816
$ws.Range("I71").Value = 2749.75
$ws.Range("J71").Value = 9666.666999999999
$ws.Range("K71").Value = 13748.75
$ws.Range("L71").Value = 48333.335
$ws.Range("M71").Value = -10004.75
$ws.Range("N71").Value = -55821.335
$ws.Range("H82").Value = 2387.1667
$ws.Range("I82").Value = 2064.6
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 2064.6
$ws.Range("L82").Value = 4000
$ws.Range("M82").Value = -1703.6
$ws.Range("N82").Value = -4722
$ws.Range("H85").Value = 2387.1667
$ws.Range("I85").Value = 2064.6
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 2064.6
$ws.Range("L85").Value = 4000
$ws.Range("M85").Value = -816.5999999999999
$ws.Range("N85").Value = -6496
$ws.Range("H93").Value = 20999.5
$ws.Range("I93").Value = 20999.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 20999.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -19751.5
$ws.Range("N93").ClearContents()
$ws.Range("H136").Value = 5224.1113
$ws.Range("I136").Value = 5224.1113
$ws.Range("K136").Value = 15672.3339
$ws.Range("M136").Value = -13122.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
